$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C10 value changes from 18 to 1 (stored as numeric)
$ws.Range("C10").Value = 1
